$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.383.50"
$ws.Range("E2").Value = "'  +0.48%  "
$ws.Range("D3").Value = "'1.690.75"
$ws.Range("E3").Value = "'  +0.08%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "'  +0.45%  "
$ws.Range("D5").Value = "'218.84"
$ws.Range("E5").Value = "'  -0.08%  "
$ws.Range("D6").Value = "'0.5471"
$ws.Range("E6").Value = "'  +4.03%  "
$ws.Range("E7").Value = "'  +0.44%  "
$ws.Range("D8").Value = "'0.2733"
$ws.Range("E8").Value = "'  +1.27%  "
$ws.Range("D9").Value = "'0.06464"
$ws.Range("E9").Value = "'  +0.26%  "
$ws.Range("D10").Value = "'22.00"
$ws.Range("E10").Value = "'  -0.31%  "
$ws.Range("D11").Value = "'0.07677"
$ws.Range("D12").Value = "'1.695.62"
$ws.Range("E12").Value = "'  +0.10%  "
$ws.Range("D13").Value = "'4.543"
$ws.Range("E13").Value = "'  -0.30%  "
$ws.Range("D14").Value = "'0.5833"
$ws.Range("E14").Value = "'  -0.51%  "
$ws.Range("D15").Value = "'0.000008395"
$ws.Range("E15").Value = "'  -1.85%  "
$ws.Range("D16").Value = "'65.32"
$ws.Range("E16").Value = "'  +1.04%  "
$ws.Range("D17").Value = "'26.436.13"
$ws.Range("E17").Value = "'  +0.49%  "
$ws.Range("D18").Value = "'4.949"
$ws.Range("E18").Value = "'  -0.43%  "
$ws.Range("E19").Value = "'  +0.43%  "
$ws.Range("D20").Value = "'10.97"
$ws.Range("E20").Value = "'  +0.95%  "
$ws.Range("D21").Value = "'191.42"
$ws.Range("E21").Value = "'  +0.45%  "
$ws.Range("D22").Value = "'6.256"
$ws.Range("E22").Value = "'  +0.13%  "
$ws.Range("E23").Value = "'  +0.48%  "
$ws.Range("D24").Value = "'149.64"
$ws.Range("E24").Value = "'  +3.17%  "
$ws.Range("D25").Value = "'0.1322"
$ws.Range("E25").Value = "'  +7.12%  "
$ws.Range("D26").Value = "'7.891"
$ws.Range("E26").Value = "'  +2.66%  "
$ws.Range("D27").Value = "'15.73"
$ws.Range("E27").Value = "'  -0.87%  "
$ws.Range("D28").Value = "'0.06334"
$ws.Range("E28").Value = "'  -5.20%  "
$ws.Range("D29").Value = "'1.407"
$ws.Range("E29").Value = "'  +3.76%  "
$ws.Range("E30").Value = "'  -0.03%  "
$ws.Range("D31").Value = "'3.606"
$ws.Range("E31").Value = "'  +0.35%  "
$ws.Range("D32").Value = "'3.594"
$ws.Range("E32").Value = "'  +1.08%  "
$ws.Range("D33").Value = "'1.683"
$ws.Range("E33").Value = "'  +0.70%  "
$ws.Range("E34").Value = "'  +1.45%  "
$ws.Range("D35").Value = "'0.6170"
$ws.Range("E35").Value = "'  -0.74%  "
$ws.Range("D36").Value = "'2.411"
$ws.Range("E36").Value = "'  +0.96%  "
$ws.Range("D37").Value = "'2.709"
$ws.Range("E37").Value = "'  -0.06%  "
$ws.Range("D38").Value = "'6.264"
$ws.Range("E38").Value = "'  -0.40%  "
$ws.Range("D39").Value = "'1.121.56"
$ws.Range("E39").Value = "'  +1.82%  "
$ws.Range("D40").Value = "'0.01635"
$ws.Range("E40").Value = "'  +0.93%  "
$ws.Range("D41").Value = "'0.8788"
$ws.Range("E41").Value = "'  -1.21%  "
$ws.Range("D42").Value = "'1.016"
$ws.Range("E42").Value = "'  +0.00%  "
$ws.Range("D43").Value = "'101.77"
$ws.Range("E43").Value = "'  +0.63%  "
$ws.Range("D44").Value = "'1.842.33"
$ws.Range("E44").Value = "'  +0.24%  "
$ws.Range("D45").Value = "'0.00000000109"
$ws.Range("E45").Value = "'  -5.60%  "
$ws.Range("D46").Value = "'57.48"
$ws.Range("E46").Value = "'  +1.06%  "
$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'8.200"
$ws.Range("E47").Value = "'  +0.13%  "
$ws.Range("B48").Value = "'Frax"
$ws.Range("C48").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.011"
$ws.Range("E48").Value = "'  +0.16%  "
$ws.Range("D49").Value = "'0.05284"
$ws.Range("E49").Value = "'  +0.44%  "
$ws.Range("D50").Value = "'6.109"
$ws.Range("E50").Value = "'  +1.02%  "
